# Applies the authoritative edit: rows 2-6 are cyclically rotated so that
# each rows full set of field values moves to a new row position
# (new row 2 <= old row 5, new row 3 <= old row 6, new row 4 <= old row 2,
#  new row 5 <= old row 3, new row 6 <= old row 4), matching the upstream
# re-export/re-ordering of the underlying Artportalen observation rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 96483446
$ws.Range("B2").Value2 = 89392
$ws.Range("E2").Value2 = 1202
$ws.Range("F2").Value2 = 'Ullticka'
$ws.Range("G2").Value2 = 'Phellinidium ferrugineofuscum'
$ws.Range("H2").Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("N2").Value2 = ''
$ws.Range("Q2").Value2 = 358350.7523457828
$ws.Range("R2").Value2 = 6853048.491440171
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value2 = '2021-10-05'
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value2 = '2021-10-05'
$ws.Range("AC2").Value2 = 'Avverkningsanmält.'
$ws.Range("AJ2").Value2 = 'gran'
$ws.Range("AK2").Value2 = 'Picea abies'
$ws.Range("AM2").Value2 = 'Liggande död trädstam, utan markontakt'
$ws.Range("AO2").Value2 = 'Horizontal, dead without ground contact # Picea abies'
$ws.Range("AW2").Value2 = 'John-Olof Halvarsson'
$ws.Range("AX2").Value2 = 'John-Olof Halvarsson'

# Row 3
$ws.Range("A3").Value2 = 96483356
$ws.Range("B3").Value2 = 77259
$ws.Range("D3").Value2 = 'NT'
$ws.Range("E3").Value2 = 228912
$ws.Range("F3").Value2 = 'Mörk kolflarnlav'
$ws.Range("G3").Value2 = 'Carbonicola myrmecina'
$ws.Range("H3").Value2 = '(Ach.) Bendiksby & Timdal'
$ws.Range("I3").Value2 = ''
$ws.Range("J3").Value2 = ''
$ws.Range("N3").Value2 = ''
$ws.Range("Q3").Value2 = 358049.7783280805
$ws.Range("R3").Value2 = 6852584.895215719
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value2 = '2021-10-05'
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value2 = '2021-10-05'
$ws.Range("AC3").Value2 = 'Avverkningsanmält.'
$ws.Range("AW3").Value2 = 'John-Olof Halvarsson'
$ws.Range("AX3").Value2 = 'John-Olof Halvarsson'

# Row 4
$ws.Range("A4").Value2 = 73741071
$ws.Range("B4").Value2 = 89410
$ws.Range("D4").Value2 = 'NT'
$ws.Range("E4").Value2 = 5432
$ws.Range("F4").Value2 = 'Granticka'
$ws.Range("G4").Value2 = 'Porodaedalea chrysoloma'
$ws.Range("H4").Value2 = '(Fr.) Fiasson & Niemelä'
$ws.Range("L4").ClearContents()
$ws.Range("Q4").Value2 = 358226.9757352364
$ws.Range("R4").Value2 = 6853003.403278765
$ws.Range("AC4").ClearContents()

# Row 5
$ws.Range("A5").Value2 = 73741042
$ws.Range("B5").Value2 = 89633
$ws.Range("D5").Value2 = 'VU'
$ws.Range("E5").Value2 = 65
$ws.Range("F5").Value2 = 'Fläckporing'
$ws.Range("G5").Value2 = 'Anthoporia albobrunnea'
$ws.Range("H5").Value2 = '(Romell) Karasiński & Niemelä'
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value2 = '1'
$ws.Range("J5").Value2 = 'mycel'
$ws.Range("N5").ClearContents()
$ws.Range("Q5").Value2 = 358011.5696047437
$ws.Range("R5").Value2 = 6852865.0116787
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value2 = '2018-10-26'
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value2 = '2018-10-26'
$ws.Range("AC5").ClearContents()
$ws.Range("AJ5").ClearContents()
$ws.Range("AK5").ClearContents()
$ws.Range("AM5").ClearContents()
$ws.Range("AO5").ClearContents()
$ws.Range("AW5").Value2 = 'Andreas Öster'
$ws.Range("AX5").Value2 = 'Andreas Öster'

# Row 6
$ws.Range("A6").Value2 = 73741060
$ws.Range("B6").Value2 = 94532
$ws.Range("D6").Value2 = 'EN'
$ws.Range("E6").Value2 = 1452
$ws.Range("F6").Value2 = 'Timmerskapania'
$ws.Range("G6").Value2 = 'Scapania apiculata'
$ws.Range("H6").Value2 = 'Spruce'
$ws.Range("L6").Value2 = ''
$ws.Range("N6").ClearContents()
$ws.Range("Q6").Value2 = 358240.2044175995
$ws.Range("R6").Value2 = 6852979.173113798
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value2 = '2018-10-26'
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value2 = '2018-10-26'
$ws.Range("AC6").Value2 = 'På tallved i stor lok, torrlagd för tillfället. Ej sågade kanter. Röda groddkorn. Lämnar kollekt till U.Gunnarsson. Granskog runt loken. Kanske totalt 3 granlågor o 1 tallåga äldre. Skulle behöva kompleteras med mer ved. Dålig kantzon från hygget.'
$ws.Range("AW6").Value2 = 'Andreas Öster'
$ws.Range("AX6").Value2 = 'Andreas Öster'
